# update DORA for better rendering
# remove ref_id for articles

$wb = $excel.ActiveWorkbook

$wsLibrary = $wb.Worksheets.Item("library_content")
$wsDora = $wb.Worksheets.Item("dora")

# Remove the Ref_id (column C) values for the article rows, keeping the
# existing cell formatting (style) intact.
$refIdRows = @(4, 23, 42, 48, 56, 73, 79, 99, 115, 128, 132, 146, 156, 165, 182, 185, 192, 196, 211, 225, 259, 269, 297)
foreach ($r in $refIdRows) {
    $wsDora.Cells.Item($r, 3).ClearContents()
}

# Update the view state for both sheets: move the current selection on
# "library_content" down to B11 and move the current selection / active
# tab to the "dora" sheet at C297.
$wsLibrary.Activate()
$wsLibrary.Range("B11").Select()

$wsDora.Activate()
$wsDora.Range("C297").Select()
